$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190641045570374
$ws.Range("B1").Value = 2.472947835922241
$ws.Range("D1").Value = 2.27983021736145
$ws.Range("E1").Value = 1.1801438331604
